$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that values such as
# "27.375.63" or "1.000" are not reinterpreted/rounded as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '27.375.63'
$ws.Range("D3").Value = '1.826.83'
$ws.Range("D5").Value = '313.15'
$ws.Range("D7").Value = '0.4462'
$ws.Range("D8").Value = '0.3768'
$ws.Range("D9").Value = '0.07416'
$ws.Range("D10").Value = '0.8783'
$ws.Range("D12").Value = '1.828.20'
$ws.Range("D13").Value = '6.715'
$ws.Range("D14").Value = '5.430'
$ws.Range("D16").Value = '0.07067'
$ws.Range("D18").Value = '0.000008819'
$ws.Range("D19").Value = '1.000'
$ws.Range("D21").Value = '27.375.86'
$ws.Range("D22").Value = '5.356'
$ws.Range("D23").Value = '10.94'
$ws.Range("D24").Value = '1.960'
$ws.Range("D25").Value = '151.28'
$ws.Range("D26").Value = '2.283'
$ws.Range("D28").Value = '5.350'
$ws.Range("D29").Value = '117.16'
$ws.Range("D30").Value = '0.08916'
$ws.Range("D31").Value = '0.7933'
$ws.Range("D32").Value = '1.199'
$ws.Range("D34").Value = '2.963'
$ws.Range("D35").Value = '0.9999'
$ws.Range("D36").Value = '1.107'
$ws.Range("D37").Value = '0.01981'
$ws.Range("D38").Value = '0.05277'
$ws.Range("D39").Value = '7.371'
$ws.Range("D40").Value = '0.5340'
$ws.Range("D41").Value = '2.877'
$ws.Range("D42").Value = '2.329'
$ws.Range("D44").Value = '8.691'
$ws.Range("D45").Value = '0.5079'
$ws.Range("D46").Value = '10.63'
$ws.Range("D47").Value = '105.39'
$ws.Range("D48").Value = '1.690'
$ws.Range("D49").Value = '0.9998'
$ws.Range("D50").Value = '0.06392'

# Restore the default "Normal" style on column D so no stray number-format
# style is left behind (matches original workbook formatting).
$ws.Range("D2:D51").Style = "Normal"

# --- Column E (Volume/1h) updates ---
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  +2.71%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("E22").Value = '  +3.80%  '
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  +3.70%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +2.69%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("E31").Value = '  +5.57%  '
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("E33").Value = '  +2.26%  '
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +1.51%  '
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  +4.85%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  +18.59%  '
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("E44").Value = '  +2.94%  '
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  +5.64%  '

